$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 - logistic_embeddings
$ws.Range("C5").Value = 0.456
$ws.Range("D5").Value = 0.566
$ws.Range("E5").Value = 0.592
$ws.Range("F5").Value = 0.64
$ws.Range("G5").Value = 0.651
$ws.Range("H5").Value = 0.663

# Row 7 - classical-best-embed (renamed)
$ws.Range("A7").Value = "classical-best-embed"
$ws.Range("C7").Value = 0.456
$ws.Range("E7").Value = 0.592
$ws.Range("F7").Value = 0.64

# Row 8 - BERT-base
$ws.Range("C8").Value = 0.462
$ws.Range("D8").Value = 0.632
$ws.Range("E8").Value = 0.661
$ws.Range("F8").Value = 0.694
$ws.Range("G8").Value = 0.719
$ws.Range("H8").Value = 0.73

# Row 9 - BERT-base-nli
$ws.Range("B9").Value = 0.314
$ws.Range("C9").Value = 0.509
$ws.Range("D9").Value = 0.645
$ws.Range("E9").Value = 0.667
$ws.Range("F9").Value = 0.681
$ws.Range("G9").Value = 0.709
$ws.Range("H9").Value = 0.719
